# Auto-generated-assisted Excel COM-interop edit script
# Implements the weekly crime-stat data refresh described in the commit message
# 'New crime data collected': updates the report header (volume number + week-of
# dates) and the Citywide crime-complaint figures for rows 14-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: bump the report volume number (4 -> 5) ---
$ws.Range("A8").Characters(21, 1).Text = "5"

# --- Header text: the week-of date range shifts forward one week ---
$ws.Range("C9").Characters(27, 9).Text = "1/30/2023"
$ws.Range("C9").Characters(47, 9).Text = "2/5/2023"

# --- Citywide crime-complaint table (rows 14-30, columns C:N) ---
# Row 14
$ws.Range("C14").Value = 5
$ws.Range("E14").Value = -16.666666666666
$ws.Range("F14").Value = 25
$ws.Range("G14").Value = 28
$ws.Range("H14").Value = -10.714285714285
$ws.Range("I14").Value = 32
$ws.Range("J14").Value = 35
$ws.Range("K14").Value = -8.571428571428
$ws.Range("L14").Value = -13.513513513513
$ws.Range("M14").Value = -27.272727272727
$ws.Range("N14").Value = -84.466019417475

# Row 15
$ws.Range("C15").Value = 20
$ws.Range("D15").Value = 33
$ws.Range("E15").Value = -39.393939393939
$ws.Range("F15").Value = 102
$ws.Range("G15").Value = 126
$ws.Range("H15").Value = -19.047619047619
$ws.Range("I15").Value = 145
$ws.Range("J15").Value = 159
$ws.Range("K15").Value = -8.80503144654
$ws.Range("L15").Value = 27.19298245614
$ws.Range("M15").Value = 21.848739495798
$ws.Range("N15").Value = -50

# Row 16
$ws.Range("C16").Value = 283
$ws.Range("D16").Value = 247
$ws.Range("E16").Value = 14.574898785425
$ws.Range("F16").Value = 1157
$ws.Range("G16").Value = 1113
$ws.Range("H16").Value = 3.953279424977
$ws.Range("I16").Value = 1525
$ws.Range("J16").Value = 1435
$ws.Range("K16").Value = 6.271777003484
$ws.Range("L16").Value = 45.793499043977
$ws.Range("M16").Value = -21.794871794871
$ws.Range("N16").Value = -82.934198746642

# Row 17
$ws.Range("C17").Value = 404
$ws.Range("D17").Value = 389
$ws.Range("E17").Value = 3.856041131105
$ws.Range("F17").Value = 1730
$ws.Range("H17").Value = 8.057464084946
$ws.Range("I17").Value = 2333
$ws.Range("J17").Value = 2079
$ws.Range("K17").Value = 12.217412217412
$ws.Range("L17").Value = 30.627099664053
$ws.Range("M17").Value = 62.352122477383
$ws.Range("N17").Value = -33.740414654927

# Row 18
$ws.Range("C18").Value = 255
$ws.Range("D18").Value = 265
$ws.Range("E18").Value = -3.77358490566
$ws.Range("F18").Value = 1099
$ws.Range("G18").Value = 1107
$ws.Range("H18").Value = -0.722673893405
$ws.Range("I18").Value = 1452
$ws.Range("J18").Value = 1414
$ws.Range("K18").Value = 2.687411598302
$ws.Range("L18").Value = 17.952883834281
$ws.Range("M18").Value = -25.538461538461
$ws.Range("N18").Value = -86.063921681543

# Row 19
$ws.Range("C19").Value = 912
$ws.Range("D19").Value = 884
$ws.Range("E19").Value = 3.167420814479
$ws.Range("F19").Value = 3624
$ws.Range("G19").Value = 3741
$ws.Range("H19").Value = -3.127506014434
$ws.Range("I19").Value = 4620
$ws.Range("J19").Value = 4762
$ws.Range("K19").Value = -2.981940361192
$ws.Range("L19").Value = 60.027710426047
$ws.Range("M19").Value = 31.25
$ws.Range("N19").Value = -41.754916792738

# Row 20
$ws.Range("C20").Value = 261
$ws.Range("D20").Value = 227
$ws.Range("E20").Value = 14.977973568281
$ws.Range("F20").Value = 1082
$ws.Range("G20").Value = 1048
$ws.Range("H20").Value = 3.24427480916
$ws.Range("I20").Value = 1397
$ws.Range("J20").Value = 1333
$ws.Range("K20").Value = 4.801200300075
$ws.Range("L20").Value = 100.143266475645
$ws.Range("M20").Value = 40.826612903225
$ws.Range("N20").Value = -88.273314866112

# Row 21
$ws.Range("C21").Value = 2140
$ws.Range("D21").Value = 2051
$ws.Range("E21").Value = 4.339346660165
$ws.Range("F21").Value = 8819
$ws.Range("G21").Value = 8764
$ws.Range("H21").Value = 0.627567320858
$ws.Range("I21").Value = 11504
$ws.Range("J21").Value = 11217
$ws.Range("K21").Value = 2.558616385842
$ws.Range("L21").Value = 47.506090524426
$ws.Range("M21").Value = 14.902117459049
$ws.Range("N21").Value = -73.380845500613

# Row 22
$ws.Range("D22").Value = 35
$ws.Range("E22").Value = -2.857142857142
$ws.Range("F22").Value = 125
$ws.Range("G22").Value = 176
$ws.Range("H22").Value = -28.977272727272
$ws.Range("I22").Value = 169
$ws.Range("J22").Value = 228
$ws.Range("K22").Value = -25.877192982456
$ws.Range("L22").Value = 32.03125
$ws.Range("M22").Value = -25.877192982456

# Row 23
$ws.Range("C23").Value = 110
$ws.Range("D23").Value = 111
$ws.Range("E23").Value = -0.9009009009
$ws.Range("F23").Value = 413
$ws.Range("G23").Value = 401
$ws.Range("H23").Value = 2.992518703241
$ws.Range("I23").Value = 552
$ws.Range("J23").Value = 520
$ws.Range("K23").Value = 6.153846153846
$ws.Range("L23").Value = 16.949152542372
$ws.Range("M23").Value = 52.908587257617

# Row 24
$ws.Range("C24").Value = 1911
$ws.Range("D24").Value = 1792
$ws.Range("E24").Value = 6.640625
$ws.Range("F24").Value = 8240
$ws.Range("G24").Value = 7416
$ws.Range("H24").Value = 11.111111111111
$ws.Range("I24").Value = 10287
$ws.Range("J24").Value = 9152
$ws.Range("K24").Value = 12.40166083916
$ws.Range("L24").Value = 37.913929481163
$ws.Range("M24").Value = 34.242463787028

# Row 25
$ws.Range("C25").Value = 735
$ws.Range("D25").Value = 729
$ws.Range("E25").Value = 0.823045267489
$ws.Range("F25").Value = 2964
$ws.Range("G25").Value = 2735
$ws.Range("H25").Value = 8.372943327239
$ws.Range("I25").Value = 3816
$ws.Range("J25").Value = 3421
$ws.Range("K25").Value = 11.546331482022
$ws.Range("L25").Value = 42.760942760942
$ws.Range("M25").Value = -4.647676161919

# Row 26
$ws.Range("C26").Value = 45
$ws.Range("D26").Value = 50
$ws.Range("E26").Value = -10
$ws.Range("F26").Value = 183
$ws.Range("G26").Value = 190
$ws.Range("H26").Value = -3.684210526315
$ws.Range("I26").Value = 237
$ws.Range("J26").Value = 235
$ws.Range("K26").Value = 0.851063829787
$ws.Range("L26").Value = 15.609756097561

# Row 27
$ws.Range("C27").Value = 86
$ws.Range("D27").Value = 86
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 354
$ws.Range("G27").Value = 295
$ws.Range("H27").Value = 20
$ws.Range("I27").Value = 462
$ws.Range("J27").Value = 375
$ws.Range("K27").Value = 23.2
$ws.Range("L27").Value = 35.882352941176

# Row 28
$ws.Range("C28").Value = 16
$ws.Range("D28").Value = 25
$ws.Range("E28").Value = -36
$ws.Range("F28").Value = 78
$ws.Range("G28").Value = 100
$ws.Range("H28").Value = -22
$ws.Range("I28").Value = 97
$ws.Range("J28").Value = 127
$ws.Range("K28").Value = -23.622047244094
$ws.Range("L28").Value = 2.105263157894
$ws.Range("M28").Value = -11.818181818181
$ws.Range("N28").Value = -83.130434782608

# Row 29
$ws.Range("C29").Value = 15
$ws.Range("D29").Value = 22
$ws.Range("E29").Value = -31.818181818181
$ws.Range("F29").Value = 65
$ws.Range("G29").Value = 88
$ws.Range("H29").Value = -26.136363636363
$ws.Range("I29").Value = 83
$ws.Range("J29").Value = 114
$ws.Range("K29").Value = -27.19298245614
$ws.Range("L29").Value = -4.597701149425
$ws.Range("M29").Value = -15.306122448979
$ws.Range("N29").Value = -84.339622641509

# Row 30
$ws.Range("C30").Value = 6
$ws.Range("E30").Value = -33.333333333333
$ws.Range("F30").Value = 27
$ws.Range("G30").Value = 38
$ws.Range("H30").Value = -28.947368421052
$ws.Range("I30").Value = 35
$ws.Range("J30").Value = 46
$ws.Range("K30").Value = -23.91304347826
$ws.Range("L30").Value = 75

